$d = $word.ActiveDocument

# Availability SLA corrected from 99.95% to 99.5% (all occurrences:
# bullet heading, summary table cell, risk/assumptions table cells,
# KPI table cell, and checklist bullet).
$d.Content.Find.Execute("99.95%", $false, $false, $false, $false, $false, `
    $true, 1, $false, "99.5%", 2)

# Radio technology corrected from GSM-R to TETRA (heading label and
# the descriptive sentence referencing the international standard).
$d.Content.Find.Execute("TETRA-GSM-R:", $false, $false, $false, $false, $false, `
    $true, 1, $false, "TETRA-TETRA:", 2)

$d.Content.Find.Execute("Comunicaciones GSM-R estándar internacional", $false, $false, `
    $false, $false, $false, $true, 1, $false, `
    "Comunicaciones TETRA estándar internacional", 2)
